# B6-PowerPoint.pptx — re-style the three data tables.
#
# The deck's three tables (slides 14, 15 and 16 — each table is the first
# shape on its slide) were re-pointed from the plain "No Style, No Grid"
# table style to a themed table style:
#   {2EE0FCEA-A1A4-41E8-8348-AA9BF06DDD66}  ->  {D780541B-B088-4466-A431-7C725DE1E136}
#
# Table.Style is read-only in the PowerPoint object model; the supported
# way to re-point a table at a different style GUID is Table.ApplyStyle.

$p = $ppt.ActivePresentation

$oldStyleId = "{2EE0FCEA-A1A4-41E8-8348-AA9BF06DDD66}"
$newStyleId = "{D780541B-B088-4466-A431-7C725DE1E136}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    for ($shapeIndex = 1; $shapeIndex -le $slide.Shapes.Count; $shapeIndex++) {
        $shape = $slide.Shapes.Item($shapeIndex)
        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.Style -eq $oldStyleId) {
                $table.ApplyStyle($newStyleId)
            }
        }
    }
}
